# Swap the values of columns A, B, D, E, F, G, H, Q, R, S between row 2 and row 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "S")

foreach ($col in $columns) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $valRow2 = $cellRow2.Value2
    $valRow3 = $cellRow3.Value2

    $cellRow2.Value = $valRow3
    $cellRow3.Value = $valRow2
}
